# JS-SPA-Self-Evaluation-Protocol.xlsx
# Commit: "edit user profile done"
#
# The "Edit User Profile" row (B29) score goes from 0 -> 3, the
# "Numbers of Commits in GitHub" row (B9) score goes from 22 -> 23, and
# the Total Score (C51, =SUM(C6:C50)) recalculates from 127 -> 131.
# The sheet's view is also left scrolled/selected on the row that was
# just edited (C29), matching where the user's cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numbers of Commits in GitHub: 22 -> 23
$ws.Range("C9").Value = 23

# Edit User Profile: 0 -> 3
$ws.Range("C29").Value = 3

# Leave the selection/active cell on C29 (where the edit happened) and
# scroll the window so row 13 is at the top, matching the saved view
# state (topLeftCell="A13", activeCell/sqref="C29").
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("C29").Select()

# C51 holds =SUM(C6:C50); it auto-recalculates to 131 once C9/C29 change.
